# Apply the Alvearie FHIR IG "eng-conversation-type" ValueSet metadata refresh:
#   - Version bump 5.0.0 -> 6.0.0
#   - Date bump to the new publish timestamp
#   - Publisher value filled in ("Alvearie Team")
#   - "Contact" / "No display for ContactDetail" row replaced by a new
#     "Jurisdiction" / "United States of America" row
#   - everything below shifts up by one row (the old trailing duplicate
#     Contact row goes away)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$ws.Range("B9").Value = "Alvearie Team"

# Replace the old "Contact" row with a new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Shift the remaining rows up by one (drop the duplicate Contact row that
# used to sit at row 11)
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "Codes indicating types for Engagement communications"

$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""

$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# Remove the now-unused trailing row 15
$ws.Rows.Item(15).Delete()
